# Swap the presentation's theme colour scheme from "Integral" to the
# stock "Office Theme" colours (ppt/theme/theme1.xml, used by
# SlideMaster1 -> every layout -> every slide).
#
# The font scheme (Arial everywhere) and format scheme (fills / lines /
# effects) are already byte-identical between the Integral and Office
# themes in this deck, so only the 12 theme colour slots - and nothing
# else - need to change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (R,G,B) for the standard Office theme palette
$officeColors = @(
    @(0,0,0),         # 1  dk1
    @(255,255,255),   # 2  lt1
    @(68,84,106),      # 3  dk2
    @(231,230,230),   # 4  lt2
    @(91,155,213),    # 5  accent1
    @(237,125,49),    # 6  accent2
    @(165,165,165),   # 7  accent3
    @(255,192,0),     # 8  accent4
    @(68,114,196),    # 9  accent5
    @(112,173,71),    # 10 accent6
    @(5,99,193),      # 11 hlink
    @(149,79,114)     # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $c = $officeColors[$i - 1]
    $r = $c[0]
    $g = $c[1]
    $b = $c[2]
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
